$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CRMAccuracyData")

# --- Corrected CRM batch values for existing titration rows (135-140) ---
$ws.Range("C135").Value = 2215.3200000000002
$ws.Range("C136").Value = 2215.3200000000002
$ws.Range("C137").Value = 2224.4699999999998
$ws.Range("C138").Value = 2224.4699999999998
$ws.Range("C139").Value = 2224.4699999999998
$ws.Range("C140").Value = 2224.4699999999998

# --- Row 141: new titration result replaces the previous duplicate entry ---
$ws.Range("B139").Copy()
$ws.Range("B141").PasteSpecial(-4122)
$ws.Range("B141").Value = 2209.71801
$ws.Range("C141").Value = 2224.4699999999998
$ws.Range("D141").Formula = "=100*(B141-C141)/C141"

# --- Row 142: brand-new titration entry (dmb training - sarah and eliza) ---
$ws.Range("B139").Copy()
$ws.Range("B142").PasteSpecial(-4122)
$ws.Range("C139").Copy()
$ws.Range("C142").PasteSpecial(-4122)
$ws.Range("D139").Copy()
$ws.Range("D142").PasteSpecial(-4122)

$ws.Range("B142").Value = 2202.70937
$ws.Range("C142").Value = 2215.3200000000002
$ws.Range("D142").Formula = "=100*(B142-C142)/C142"

# --- Restore selection state to match the saved view ---
$ws.Range("D141").Select()
